$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New data rows (8-11): "Book Moves" as Player2
# Set B8 first so shared string "Book Moves" is minted before "P1Winrate"
$ws.Cells.Item(8, 2).Value = "Book Moves"

# New header for column F
$ws.Range("F1").Value = "P1Winrate"

$bookRows = @(
    @(106, 90, 4),
    @(29, 21, 0),
    @(94, 55, 1),
    @(102, 44, 4)
)

$r = 8
foreach ($row in $bookRows) {
    $ws.Cells.Item($r, 1).Value = "Vanilla MCTS"
    $ws.Cells.Item($r, 2).Value = "Book Moves"
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $r++
}

# Row 12: Greedy Book vs Vanilla MCTS
$ws.Cells.Item(12, 1).Value = "Greedy Book"
$ws.Cells.Item(12, 2).Value = "Vanilla MCTS"
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0

# F2 formula - plain formula
$ws.Range("F2").Formula = "=C2/(C2+D2)"

# F3:F12 - shared formula
$ws.Range("F3:F12").Formula = "=C3/(C3+D3)"

$ws.Range("D13").Select()
